# Auto-generated edit script: apply hybrid bold+color highlighting to quantitative metrics
$d = $word.ActiveDocument
$boldColor = 5258796  # RGB(0x2C,0x3E,0x50) in BGR-packed VBA color form -> renders as w:color val 2C3E50

$edits = @(
    @{
        Old = "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%"
        Segs = @(
        @{ Text = "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from "; Bold = $false },
        @{ Text = "23%"; Bold = $true },
        @{ Text = " to "; Bold = $false },
        @{ Text = "64%"; Bold = $true }
        )
    },
    @{
        Old = "• Utilized advanced sampling methods to decrease survey margin of error from ±4.2% to ±2.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes"
        Segs = @(
        @{ Text = "• Utilized advanced sampling methods to decrease survey margin of error from "; Bold = $false },
        @{ Text = "±4.2%"; Bold = $true },
        @{ Text = " to "; Bold = $false },
        @{ Text = "±2.1%"; Bold = $true },
        @{ Text = ", increasing voter turnout prediction accuracy from "; Bold = $false },
        @{ Text = "71%"; Bold = $true },
        @{ Text = " to "; Bold = $false },
        @{ Text = "87%"; Bold = $true },
        @{ Text = ", and ensuring survey results more closely reflected true population attitudes"; Bold = $false }
        )
    },
    @{
        Old = "• Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis"
        Segs = @(
        @{ Text = "• Trigonometric algorithm for boundary estimation reduced mapping costs by "; Bold = $false },
        @{ Text = "73.5%"; Bold = $true },
        @{ Text = ", saving campaigns and organizations "; Bold = $false },
        @{ Text = "`$4.7M"; Bold = $true },
        @{ Text = " and enabling smaller nonprofits to conduct analysis"; Bold = $false }
        )
    },
    @{
        Old = "• Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion"
        Segs = @(
        @{ Text = "• Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over "; Bold = $false },
        @{ Text = "`$2"; Bold = $true },
        @{ Text = " trillion"; Bold = $false }
        )
    },
    @{
        Old = "• Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%"
        Segs = @(
        @{ Text = "• Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by "; Bold = $false },
        @{ Text = "57%"; Bold = $true }
        )
    },
    @{
        Old = "• Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from ±4.2% to ±2.1%"
        Segs = @(
        @{ Text = "• Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from "; Bold = $false },
        @{ Text = "±4.2%"; Bold = $true },
        @{ Text = " to "; Bold = $false },
        @{ Text = "±2.1%"; Bold = $true }
        )
    },
    @{
        Old = "• Increased voter turnout prediction accuracy from 71% to 87%"
        Segs = @(
        @{ Text = "• Increased voter turnout prediction accuracy from "; Bold = $false },
        @{ Text = "71%"; Bold = $true },
        @{ Text = " to "; Bold = $false },
        @{ Text = "87%"; Bold = $true }
        )
    },
    @{
        Old = "• Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%"
        Segs = @(
        @{ Text = "• Methodological advancement: Improved segmentation accuracy "; Bold = $false },
        @{ Text = "34%"; Bold = $true },
        @{ Text = " and survey incidence "; Bold = $false },
        @{ Text = "28%"; Bold = $true }
        )
    },
)

$appliedCount = 0
$missingCount = 0

foreach ($edit in $edits) {
    $target = $edit.Old
    $targetPara = $null
    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text
        if ($t.Length -ge 1) {
            $tt = $t.Substring(0, $t.Length - 1)
        } else {
            $tt = $t
        }
        if ($tt -eq $target) {
            $targetPara = $p
        }
    }
    if ($null -eq $targetPara) {
        Write-Output "WARNING: paragraph not found for: $target"
        $missingCount = $missingCount + 1
        continue
    }
    # Walk the paragraph left-to-right, re-slicing it into plain / bold+colored
    # runs using absolute document offsets. Only formatting is touched (no
    # text insert/delete), so offsets computed up front stay valid as we go.
    $cur = $targetPara.Range.Start
    foreach ($seg in $edit.Segs) {
        $segText = $seg.Text
        $segStart = $cur
        $segEnd = $cur + $segText.Length
        if ($seg.Bold) {
            $segRange = $d.Range($segStart, $segEnd)
            $segRange.Font.Bold = 1
            $segRange.Font.Color = $boldColor
        }
        $cur = $segEnd
    }
    $appliedCount = $appliedCount + 1
}

Write-Output "Applied $appliedCount of $($edits.Count) metric-highlight edits ($missingCount missing)."